$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update hours for 2022-02-02 (row 116) from 1 to 2
$ws.Cells.Item(116, 2).Value = 2

# 2. Insert 4 rows after row 116 to make room for the two new entries and
#    reproduce the same row-gap pattern used throughout the sheet (the
#    summary block ends up at rows 123-125 instead of 119-121).
$ws.Rows("117:120").Insert()

# Remove the two extra (now blank) rows in the middle of the gap so they
# don't show up as empty rows in the saved file - matches rows 119/120
# being absent in the target sheet.
$ws.Rows("119:120").Clear()

# 3. Fill in the new entry for 2022-02-03 (row 117)
$ws.Cells.Item(117, 1).Value = 44595
$ws.Cells.Item(117, 2).Value = 3
$ws.Cells.Item(117, 3).Value = "Notification tests, redux fake store, store Provider wrapper"
$ws.Cells.Item(117, 4).Value = "client"

# 4. Fill in the second new entry, still on 2022-02-03 (row 118, no date cell)
$ws.Cells.Item(118, 1).Clear()
$ws.Cells.Item(118, 2).Value = 1
$ws.Cells.Item(118, 3).Value = "Loginbar tests"
$ws.Cells.Item(118, 4).Value = "client"

# 5. Fix up the summary rows, now at 123-125, so their formulas reference
#    the new last data row (118) instead of the old one (116).
$ws.Cells.Item(123, 2).Formula = "=SUM(B2:B118)"
$ws.Cells.Item(125, 2).Formula = "=B123/B124*100"

# 6. Update the view state to match where the author left the selection.
$ws.Application.ActiveWindow.ScrollRow = 106
$ws.Range("B116").Select()
